$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 289.3202873889254
$ws.Range("G2").Value = 6.755560187139269
$ws.Range("H2").Value = 562.3507975863433
$ws.Range("I2").Value = 0.5116800344544622
$ws.Range("J2").Value = 0.01713649116353724
$ws.Range("K2").Value = 1.122755843714779
$ws.Range("L2").Value = 0.1961882789188312
$ws.Range("M2").Value = 0.006581945548015016
$ws.Range("N2").Value = 0.4034967452089521

# Row 3
$ws.Range("F3").Value = 0.004842899754225591
$ws.Range("G3").Value = 0.001791255137212199
$ws.Range("H3").Value = 0.007871934362300815
$ws.Range("I3").Value = 0.004494392975539166
$ws.Range("J3").Value = 0.001658810752945294
$ws.Range("K3").Value = 0.007306749978235069
$ws.Range("L3").Value = 0.004939439022672729
$ws.Range("M3").Value = 0.001877944278049141
$ws.Range("N3").Value = 0.007979769487075555

# Row 4
$ws.Range("F4").Value = 289.3251302886796
$ws.Range("G4").Value = 6.75735144227648
$ws.Range("H4").Value = 562.3586695207057
$ws.Range("I4").Value = 0.5161744274300013
$ws.Range("J4").Value = 0.01879530191648253
$ws.Range("K4").Value = 1.130062593693014
$ws.Range("L4").Value = 0.2011277179415039
$ws.Range("M4").Value = 0.008459889826064156
$ws.Range("N4").Value = 0.4114765146960276
